$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("A12").Value = 131257520
$ws.Range("B12").Value = 79245
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = 'Garnlav'
$ws.Range("G12").Value = 'Alectoria sarmentosa'
$ws.Range("H12").Value = '(Ach.) Ach.'
$ws.Range("M12").ClearContents()
$ws.Range("Q12").Value = 488939
$ws.Range("R12").Value = 6665149
$ws.Range("Z12").Value = '11:41'
$ws.Range("AB12").Value = '11:41'
$ws.Range("AC12").Value = 'Gran'

# Row 13
$ws.Range("A13").Value = 131260641
$ws.Range("Q13").Value = 488859
$ws.Range("R13").Value = 6665292
$ws.Range("Z13").Value = '15:34'
$ws.Range("AB13").Value = '15:34'
$ws.Range("AC13").Value = 'Ringhack på gran.'

# Row 14
$ws.Range("A14").Value = 131257290
$ws.Range("B14").Value = 57884
$ws.Range("E14").Value = 100109
$ws.Range("F14").Value = 'Tretåig hackspett'
$ws.Range("G14").Value = 'Picoides tridactylus'
$ws.Range("H14").Value = '(Linnaeus, 1758)'
$ws.Range("M14").Value = 'äldre spår'
$ws.Range("Q14").Value = 488842
$ws.Range("R14").Value = 6665224
$ws.Range("Z14").Value = '11:26'
$ws.Range("AB14").Value = '11:26'
$ws.Range("AC14").Value = 'Ringhack på tall.'

# Row 15
$ws.Range("A15").Value = 131256673
$ws.Range("Q15").Value = 488652
$ws.Range("R15").Value = 6665282
$ws.Range("Z15").Value = '10:54'
$ws.Range("AB15").Value = '10:54'
$ws.Range("AC15").Value = 'Ringhack på tall.'

# Row 29
$ws.Range("A29").Value = 131258531
$ws.Range("Q29").Value = 488725
$ws.Range("R29").Value = 6665212
$ws.Range("Z29").Value = '13:02'
$ws.Range("AB29").Value = '13:02'
$ws.Range("AC29").Value = 'Gran'

# Row 30
$ws.Range("A30").Value = 131257239
$ws.Range("B30").Value = 57884
$ws.Range("E30").Value = 100109
$ws.Range("F30").Value = 'Tretåig hackspett'
$ws.Range("G30").Value = 'Picoides tridactylus'
$ws.Range("H30").Value = '(Linnaeus, 1758)'
$ws.Range("M30").Value = 'färska spår'
$ws.Range("Q30").Value = 488852
$ws.Range("R30").Value = 6665286
$ws.Range("Z30").Value = '11:23'
$ws.Range("AB30").Value = '11:23'
$ws.Range("AC30").Value = 'Barkfläk, hagelsalva.'

# Row 31
$ws.Range("A31").Value = 131255910
$ws.Range("B31").Value = 79245
$ws.Range("E31").Value = 6425
$ws.Range("F31").Value = 'Garnlav'
$ws.Range("G31").Value = 'Alectoria sarmentosa'
$ws.Range("H31").Value = '(Ach.) Ach.'
$ws.Range("M31").ClearContents()
$ws.Range("Q31").Value = 488763
$ws.Range("R31").Value = 6665157
$ws.Range("Z31").Value = '10:03'
$ws.Range("AB31").Value = '10:03'
$ws.Range("AC31").Value = 'Tall.'

# Row 36
$ws.Range("A36").Value = 131257385
$ws.Range("B36").Value = 91830
$ws.Range("E36").Value = 5432
$ws.Range("F36").Value = 'Granticka'
$ws.Range("G36").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H36").ClearContents()
$ws.Range("Q36").Value = 488876
$ws.Range("R36").Value = 6665194
$ws.Range("Z36").Value = '11:31'
$ws.Range("AB36").Value = '11:31'
$ws.Range("AC36").Value = 'Lågstubbe.'

# Row 37
$ws.Range("A37").Value = 131260531
$ws.Range("B37").Value = 79245
$ws.Range("E37").Value = 6425
$ws.Range("F37").Value = 'Garnlav'
$ws.Range("G37").Value = 'Alectoria sarmentosa'
$ws.Range("H37").Value = '(Ach.) Ach.'
$ws.Range("Q37").Value = 488786
$ws.Range("R37").Value = 6665188
$ws.Range("Z37").Value = '15:25'
$ws.Range("AB37").Value = '15:25'
$ws.Range("AC37").Value = 'Gran'
